$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.006.78"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.886.11"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5160"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3746"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07184"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8998"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07645"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "1.889.07"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.237"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008486"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "27.040.95"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.059"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "2.120.44"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.378"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.906"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.784"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09199"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05043"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.227"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7641"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.991"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.274"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.583"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5602"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.070"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.095"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.634"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1500"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4815"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.53%  "
